$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "Augusto Cezar da Silva"

$ws.Range("A15").Select()
